$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.789879322052002
$ws.Range("B1").Value = 4.213298797607422
$ws.Range("C1").Value = 1.82880175113678
$ws.Range("D1").Value = 0.8784728646278381
$ws.Range("E1").Value = 0.4740420579910278
